# Applies the commit "Update gh-pages to output generated at 456a3b4" to
# the 苏州-漫展信息 workbook.
#
# Summary of the change:
#   * Sheet "展览" (sheet1): the event "昆山·2024随机宅舞&正反派对决（免费活动）"
#     (row 5) was dropped from the scrape; all following rows shift up by
#     one, and a handful of "想去人数" (want-to-go count) values were
#     refreshed to more current numbers.
#   * Sheet "演出" (sheet2): the event "苏州·小鸳鸯上清联欢会——「遇见平江」配音
#     演员专场见面会" (row 2) was dropped; following rows shift up by one.
#   * Sheet "全部类型" (sheet4), which mirrors the union of sheet1+sheet2,
#     receives the same two row removals and the same refreshed counts.
#
# Since Rows.Delete() already shifts the text content of B..I for us, we
# only need to (1) delete the right rows, (2) poke the updated F values
# into their new positions, and (3) renumber column A (the plain 0-based
# sequence number) back into a contiguous run.

$wb = $excel.ActiveWorkbook

function Renumber-ColumnA {
    param($ws, [int]$firstDataRow, [int]$lastDataRow)
    $n = 0
    for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $n
        $n = $n + 1
    }
}

# ---------------------------------------------------------------------
# Sheet1 "展览"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Remove the row for 昆山·2024随机宅舞&正反派对决（免费活动）
$ws1.Rows.Item(5).Delete()

# Refresh "想去人数" (column F) values that changed after the row shift
$sheet1Updates = @{
    3  = 557
    4  = 1588
    6  = 201
    7  = 777
    8  = 1063
    10 = 372
    11 = 67
    12 = 527
    13 = 35
    14 = 6592
    15 = 33
    18 = 172
    20 = 1036
    21 = 15784
    23 = 24
    24 = 307
    27 = 11179
    28 = 803
    29 = 4387
    33 = 314
    34 = 134
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

Renumber-ColumnA -ws $ws1 -firstDataRow 2 -lastDataRow 34

# ---------------------------------------------------------------------
# Sheet2 "演出"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

# Remove the row for 苏州·小鸳鸯上清联欢会——「遇见平江」配音演员专场见面会
$ws2.Rows.Item(2).Delete()

Renumber-ColumnA -ws $ws2 -firstDataRow 2 -lastDataRow 5

# ---------------------------------------------------------------------
# Sheet4 "全部类型" (union of sheet1 + sheet2, kept in sync manually)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

# Remove rows for the same two dropped events. Delete the higher row
# index first so the second delete's row number stays valid.
$ws4.Rows.Item(6).Delete()   # 苏州·小鸳鸯上清联欢会...
$ws4.Rows.Item(5).Delete()   # 昆山·2024随机宅舞&正反派对决（免费活动）

$sheet4Updates = @{
    3  = 557
    4  = 1588
    6  = 201
    7  = 777
    9  = 1063
    11 = 372
    12 = 67
    13 = 527
    15 = 35
    16 = 6592
    17 = 33
    20 = 172
    23 = 1036
    24 = 15784
    26 = 24
    27 = 307
    31 = 11179
    32 = 803
    33 = 4387
    37 = 314
    38 = 134
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}

Renumber-ColumnA -ws $ws4 -firstDataRow 2 -lastDataRow 38

Write-Host "Done applying edits."
